$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.889.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.15%  "

# Row 3
$ws.Range("D3").Value = "'1.635.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.82%  "

# Row 5
$ws.Range("D5").Value = "'215.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6
$ws.Range("D6").Value = "'0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("E7").Value = "  +0.80%  "

# Row 8
$ws.Range("D8").Value = "'28.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

# Row 9
$ws.Range("E9").Value = "  +0.32%  "

# Row 10
$ws.Range("D10").Value = "'0.0609"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.38%  "

# Row 11
$ws.Range("D11").Value = "'0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.12%  "

# Row 12
$ws.Range("D12").Value = "'1.869.58"
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'1.642.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "

# Row 14
$ws.Range("D14").Value = "'0.584"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.44%  "

# Row 15
$ws.Range("D15").Value = "'9.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.71%  "

# Row 16
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$ws.Range("D17").Value = "'29.906.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "

# Row 18
$ws.Range("D18").Value = "'64.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19
$ws.Range("D19").Value = "'240.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0702"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "

# Row 21
$ws.Range("E21").Value = "  +0.69%  "

# Row 22
$ws.Range("D22").Value = "'9.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.16%  "

# Row 23
$ws.Range("E23").Value = "  +1.03%  "

# Row 24
$ws.Range("D24").Value = "'2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.02%  "

# Row 25
$ws.Range("D25").Value = "'157.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "

# Row 26
$ws.Range("D26").Value = "'15.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "

# Row 27
$ws.Range("E27").Value = "  -0.64%  "

# Row 28
$ws.Range("D28").Value = "'6.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "

# Row 29
$ws.Range("E29").Value = "  +0.72%  "

# Row 30
$ws.Range("D30").Value = "'0.0488"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "

# Row 31
$ws.Range("E31").Value = "  -1.42%  "

# Row 32
$ws.Range("D32").Value = "'3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "

# Row 33
$ws.Range("D33").Value = "'3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "

# Row 34
$ws.Range("D34").Value = "'1.421.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "

# Row 35
$ws.Range("E35").Value = "  +3.27%  "

# Row 36
$ws.Range("E36").Value = "  -1.87%  "

# Row 37
$ws.Range("D37").Value = "'2.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.01%  "

# Row 38
$ws.Range("E38").Value = "  +1.18%  "

# Row 39
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("D40").Value = "'76.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.20%  "

# Row 41
$ws.Range("D41").Value = "'0.558"
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.832"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

# Row 45
$ws.Range("E45").Value = "  +0.87%  "

# Row 46
$ws.Range("E46").Value = "  -1.63%  "

# Row 47
$ws.Range("D47").Value = "'1.777.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "

# Row 48
$ws.Range("E48").Value = "  -1.85%  "

# Row 49
$ws.Range("D49").Value = "'48.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.01%  "

# Row 50
$ws.Range("D50").Value = "'92.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.16%  "

# Row 51
$ws.Range("E51").Value = "  +6.86%  "
